$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(29, 1).Value = 45873
$ws.Cells.Item(29, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 2).Value = 6329.93994140625
$ws.Cells.Item(29, 3).Value = 6330.68994140625
$ws.Cells.Item(29, 4).Value = 6271.7099609375
$ws.Cells.Item(29, 5).Value = 6271.7099609375
$ws.Cells.Item(29, 6).Value = 4842580000
$ws.Cells.Item(29, 7).Value = 0.014737100330916
$ws.Cells.Item(29, 8).Value = 1
$ws.Cells.Item(29, 9).Value = 'nintendo break sale record did delayed gta release save switch fedex corporation fdxs ceo is a hitter say jim cramer inside gamesxs ai revolution can india become the next global gamechanger intel corporation intc taught me not to trust ai chatbots say jim cramer ipo stock of the week china stock leader atour test buy point amid market fall the procter gamble company pg it just not great say jim cramer espn obtains redzone other medium asset from nfl in blockbuster deal per report this week in ai key to agentic ai adoption plus ai action plan from u china mcdonalds plan major ai expansion to streamline customer order executive say paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame oil major beat back billion fine in kazakhstan paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame the outdated management mindset holding your business back paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame i worked in the building where people were killed here how to make your building more secure paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame columbia is selling a stylish quilted jacket for and reviewer say it is wellmade paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame top action comedy anime series in paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame already hit your k limit in here what to do next cadence design system inc cdnss number show nvidia is in great shape in china say jim cramer alnylams amvuttra make powerful entrance eye firstline dominance president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report nvidia corporation nvda could be in great shape in china say jim cramer the boeing company ba a horse it ceo ha done a great job say jim cramer doj identifies tariff fraud a enforcement priority waste management inc wms ceo is doing a remarkable job say jim cramer merck co inc mrk need to move on from talking about it animal division say jim cramer starbucks corporation sbux dont worry about luckin coffee say jim cramer sarah jessica parker announces emotional and poetic end to and just like that carrier global corporation carr delivered a good quarter say jim cramer iphone pro launch in september new colour a pro chip expected india price and more otis worldwide corporation otis had a tough quarter china wa horrible say jim cramer a apple watch series is just one click away at stacksocial celestica inc cl is up huge say jim cramer ai evolution prompt sense of urgency among cfo oracle jury say tesla wa partly to blame for fatal crash trump issue nuclear threat to russiaand wall street fear gauge go vertical crypto founder take note these sec roundtable could decide the rule that govern you sugar price pressured a brazil ramp up sugar production this week in stablecoins building nextgen rail for enterprise finance tesla must pay million in damage after fatal autopilot crash jury say cnbc miami jury find tesla partly liable in autopilot crash award million in punitive damage nbc news bond etf could catch a tailwind from job data tesla tsla is found liable in fatal autopilot crash ha to pay million the latest u stock market and global trade partner react to trump new tariff no one can be that wrong trump fire labor statistic bos after weak revised job data the green revolution anand chandra cofounder executive director aryaag on fixing india agri gap with tech cocoa price sharply lower on the outlook for adequate supply what ghislaine maxwell move to a lowersecurity federal prison camp could signal reddit want to replace google a the internet top search engine investor are all in ct bond commission approves m to lower electric bill republican say it not enough coffee price plunge a supply concern ease the latest business and u trading partner worldwide respond to trump new tariff the latest business and u trading partner worldwide respond to trump new tariff nextnav and oscilloquartz partner to advance g based timing solution for critical infrastructure resmed analyst increase their forecast after strong q earnings hedge fund and insider trading news william heard warren buffett ray dalio bill ackman rokos capital management alden global capital clear channel outdoor holding inc cco and more here are reason we are not going to buy stock in friday nosediving market this is why coinbase stock coin is cratering today tipranks phoenix az lemonade stand squeeze success expands into western state with over million in sale top cryptos to join now for that could deliver big return a new bowling center wa big news in west valley bowl face endofmonth deadline corporation for public broadcasting will shutter after trump funding cut the best trader in congress how to track politician stock trade a dallas bank is confident it loss will pay off best gold stock to watch in midtiers and junior beat the major online creator led by mrbeast and mark rober want to raise million for clean water access online creator led by mrbeast and mark rober want to raise million for clean water access online creator led by mrbeast and mark rober want to raise million for clean water access former chairman of savannah economic development dy jury order tesla to pay m in autopilot crash case opening it up to other costly lawsuit jury order tesla to pay m in autopilot crash case opening it up to other costly lawsuit online creator led by mrbeast and mark rober want to raise million for clean water access tesla found partially liable for fatal crash involving driverassistance technology jury order million in damage the washington post india secret weapon for defence sovereignty through ai and zerochina tech trump injects a new dose of uncertainty in tariff a he push start date back to aug trump injects a new dose of uncertainty in tariff a he push start date back to aug these analyst boost their forecast on microstrategy following q result trump injects a new dose of uncertainty in tariff a he push start date back to aug unit moving and portable storage return a primary sponsor for timmy hill at watkins glen nascar truck series race tech stock lead selloff a soft labor data fuel fed cut bet the eye of the hurricane why the u job market ha soured economist say here what could get more expensive from trump massive tariff hike here what could get more expensive from trump massive tariff hike trump seek to leave his goldplated stain on the white house beatrice police department hosting national night out nicole leone former executive director erie humane society ha september trial date president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report north carolina u senate race turn into a battle for the middle class sarah jessica parker announces emotional and poetic end to and just like that ghislaine maxwell join famous prisoner in move to texas prison camp wave of fake credential spark political fallout in spain gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade corporation for public broadcasting which fund pb npr to close after federal aid cut nyc subway nightmare flood trap commuter in filth and failure the latest u stock market and global trade partner react to trump new tariff what ghislaine maxwell move to a lowersecurity federal prison camp could signal what new public comment rule mean for speaking at local meeting in bismarckmandan area corporation for public broadcasting will shutter after trump funding cut congressman view here the truth about impact of beautiful bill on st louis county trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev state college rally planned a another local health system cut youth trans care more american couple are turning to italy dolce vita in a quest for memorable wedding city say thursday odor came from construction of city new industrial water plant trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia cult of cruelty ghastly photo video of emaciated israeli hostage released by gaza terrorist here what could get more expensive from trump massive tariff hike red state move to slash medicaid funding for planned parenthood trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev guest columnist we must empower support girl and boy trump order nuclear submarine repositioned after provocative statement from russia editorial maxwell talk giving more life to a conspiracy shribman who a hero in the age of trump brown fulfilling the promise of the ada a trophy hunter killed a lion in zimbabwe that wa part of a research project sparking anger pensacola man traveled to orlando stealing smokey bear sign along the way official say kamala harris insufferably cackle over the government spying on american through wireless earpods western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens trump say he deploying two nuclear sub after provocative comment from exrussian leader democrat are struggling to block trump texas power grab the healing power of sunlight office cleaner who came facetoface with the manhattan gunman detail his deadly rampage corporation for public broadcasting funder of npr and pb say it will end operation within month after federal budget cut dozen of protester swarm chuck schumer kirsten gillibrands nyc office to let gaza live fourth circuit rule unanimously in favor of employee in two case involving employer allegation of breach of restrictive covenant city county official talk about animal control a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation trump plan to build million white house party room america rising part a federal judge just called out the trump administration for lying to the supreme court k nieth to retire fbi official reportedly redacted reference to trump from epstein file trump positioning two nuclear submarine near russia following russian official threat of war how immigration legislation is dangerous to immigrant judge pause removal of k migrant from honduras nicaragua and nepal in overthetop ruling democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan trump weighs in on american eagle sydney sweeney ad shiba inu and cardano set to breakout a one new token set sight on gain top cryptos to buy in trimiq weight loss official launched in uk and ireland with green tea probiotic complex free spin no deposit bonus casino online real money casino no deposit by wild casino the shock job report set off this recession alert and hold fresh clue that ai may be boosting unemployment jpmorgan say fortune crypto presale that could make you a millionaire in dogecoin ada xrp cant match ozak ai potential arialief for neuropathy supplement official launch of a new hope for sciatic nerve discomfort relief the ig choice ha gone off the rail reader commentary what to know if youre at risk of having your wage garnished over student loan debt m raised and climbing ruvi ai ruvi called ripple xrp a coinmarketcap listing made it token one of the highest in demand why parent are suing snapchat over fentanyl death tea encouraged it user to spill then the apps data got leaked npr federal reserve governor resigns creating vacancy for trump federal reserve governor resigns creating vacancy for trump upicom bragar eagel squire pc is reminds investor of class action against neogen alto bragar eagel squire pc is reminds investor of class action against neogen alto rxsight and replimune and encourages investor to contact the firm audizen drop officially launched natural tinnitus support inspired by audifort europe ev sale report tesla ha plus win in a ev share market matsato knife official launch in australia usa redefining precision and power in the kitchen warren buffett dividend stock you can buy and hold forever crushed q earnings luka doncic agrees to a threeyear million contract extension through with the la lakers the financial swamp of banking regulator who decide which american get their bank account shut down sam altman compare openai to the manhattan projectand he not joking about the risk the democrat have a mamdani problem editorial toncoin at ondo cross while cold wallet entry crypto reward outclass both weekend gaming best ipad game to try in little pepe lilpepe price prediction for show strong potential to beat solana sol and cardano ada sluggish summer in sin city la vega see big drop in tourist ap business summarybrief at pm edt ap business summarybrief at pm edt from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price texas panel advance redrawn congressional map that could take democratic seat your money exit a business or equity position without regret pepes chart turn positive litecoin gather whale support blockdag expands web reach with devs a day workweek isnt good just for team moral it also good for the bottom line teleflex raise outlook following strong q and biotronik acquisition tucson trolley tour operator seek to expand reach former trump prosecutor jack smith probed for alleged election interference u hiring slowed sharply over the summer the wall street journal mstr news michael saylor explains why strategy strc preferred stock is the firm iphone moment coindesk harvard say it will comply with dhs request to hand over employment form missed binance coin bnb expert say ruvi ai ruvi audited token is next millionaire maker a coinmarketcap listing spark high demand what to know if youre at risk of having your wage garnished over student loan debt how to safely download and install game on your laptop shib find support ondo hold while cold wallet price gas rebate make it a top crypto to watch little pepe lilpepe gather momentum with gain ahead while cardano ada nears run the trump administration take a very orwellian turn how this week avalanche of news from washington to wall street kept investor guessing westchester capital management inc trim stake in valero energy corporation nysevlo united super pty ltd in it capacity a trustee for the construction building union superannuation fund make new investment in valero energy corporation nysevlo decadeold bitcoin wallet emptied million in btc moved amid market correction how florida quietly surpassed california in solar growth ukrainian longrange drone hit multiple russian refinery factory overnight journointerior designer jim acosta guarantee trump wh ballroom is going to look like st americu appoints new community relation manager mitsubishi ufj trust banking corp purchase share of ge aerospace nysege j safra sarasin holding ag purchase share of thermo fisher scientific inc nysetmo mitsubishi ufj trust banking corp sell share of thermo fisher scientific inc nysetmo banco santander sa ha million stake in thermo fisher scientific inc nysetmo the walt disney company nysedis share bought by bank of nova scotia trust co the walt disney company nysedis stake lifted by r squared ltd the walt disney company nysedis share acquired by j safra sarasin holding ag the walt disney company nysedis share bought by j safra sarasin holding ag banco santander sa increase stock holding in the walt disney company nysedis this ai is set to flip every major altcoinare you in figma ipo surprise winner is a charity with million sharesand a famous backstory that sparked a bitter feud over an oil fortune decade ago fortune the hidden signal in oil market moor cabot inc ha million stock holding in the goldman sachs group inc nysegs moor cabot inc buy share of the goldman sachs group inc nysegs hartford investment management co cut stock position in zoetis inc nysezts how florida quietly surpassed california in solar growth cnbc lesson from a decade of empowering woman entrepreneur to lead boldly best trucking bookkeeping service three stock to consider today blockdags no vesting pas with roi end soon eth eye jup drop just of smallbusiness owner say the child care crisis is a critical problem here why pepeto price forecast will this presale meme coin surpass dogecoin shiba inu and ethereum editorial beyond the number study into subminimum must consider all ramification thing to know about palantir pltr before it report q earnings best trucking bookkeeping service it super disrespectful wnba player exasperated after sex toy thrown onto court for second time in a week claude cutoff anthropic revoke openais api access ahead of gpt launch malaga cove capital llc sell share of abbott laboratory nyseabt bayforest capital ltd ha stake in fidelity national financial inc nysefnf vestcor inc buy share of fidelity national financial inc nysefnf modera wealth management llc sell share of arista network inc nyseanet bank of nova scotia sell share of procter gamble company the nysepg bank of nova scotia increase stock holding in arista network inc nyseanet aviso wealth management acquires share of procter gamble company the nysepg aviso wealth management cut position in the boeing company nyseba j safra sarasin holding ag purchase share of abbott laboratory nyseabt pathway financial adviser llc purchase share of eli lilly and company nyselly aviso wealth management lower stock position in the boeing company nyseba procter gamble company the nysepg share bought by j safra sarasin holding ag what are the best dividend growth stock i should invest in right now hims alert bragar eagel squire pc announces that a class action lawsuit ha been steve witkoff meet with hostage family in tel aviv susan estrich the harvard fight money v principle jet holiday ad audio used in ice deportation video shes so cringe christine pelosi tease that shes running but not the way youd think video from the dream of dr king to the nightmare that is donald trump opinion from the dream of dr king to the nightmare that is donald trump opinion from the dream of dr king to the nightmare that is donald trump opinion guam d printing facility to enhance navy maintenance capability the ig choice ha gone off the rail reader commentary trump motif can be opaque i asked a feminist economist to bare them india stand firm on russian oil import despite trump sanction threat pedaling dei cycling weekly demand bike community include fat black woman or theyre racist the u military just launched it largest air exercise in year pennsylvania lawmaker say child care is broken but dont agree on how to fix it yet twice told tale city council attends team building workshop if trump save maxwell he will destroy everything he promised pope thrill hundred of thousand of young catholic at holy year youth festival pope thrill hundred of thousand of young catholic at holy year youth festival are republican doing enough to confirm trump nominee conservative critic say no texas panel advance redrawn congressional map that could take democratic seat trump wa asked about his press secretary performance his answer took a weird turn a i see it corvallis ha done nothing to help gazans how top democrat are already gearing up for online texas is moving forward with it radical redistricting plan america trucker need a place to park george mason university president keep job get a raise despite trump admins probe in which way is the bls biased donna brazile say nprpbsdefunding maga cannot compete with idea cue a hard trip over hillary right degree wrong time new graduate face discouraging odds of landing a job trump administration move to ban abortion care for veteran unspeakably cruel manhunt intensifies for suspect in montana bar mass shooting texas panel advance redrawn congressional map that could take democratic seat bill maher know trump losing credibility when these supporter smell bullst kevin oleary slam trump for whacking bls chief after disappointing job report you dont shoot the messenger trump administration halted civil right lawsuit targeting abuse of prisoner and the mentally ill it wa never about safety one week in here whats blocked by the uk orwellian online bill these three country dominate the fighter jet export marketfor now tiktok can shape america next generation and beijing know it former trump prosecutor jack smith probed for alleged election interference virginia end fiscal year with million surplus for next budget bonus editorial cartoon for aug other voice conserving enhancing agricultural resource with solar energy trump alum share chilling google message from before secondterm return lawfare at it finest texas house panel advance redrawn congressional map that would add more gop seat letter nd representative helping put our nation closer to fiscal crisis tracey weiss our ocean backyard coming soon opening the gate to cotonicoast dairy senate make progress in averting government shutdown much earlier than usual no piecemeal deal witkoff tell hostage family trump want full gaza agreement will the wnba approve record m deal to bring connecticut sun to boston proud democrat call into chris cuomos show admits she wa wrong about trump shes not the only one texas house committee advance gopfriendly map texas house committee advance gopfriendly map texas house committee advance gopfriendly map release the epstein file and end the coverup texas house committee advance gopfriendly map texas house committee advance gopfriendly map texas house committee advance gopfriendly map the trump administration take a very orwellian turn new photo reveal iconic white house rose garden paved over after trump makeover the nfls new kickoff rule is here to stay with a slight tweak what else is changing oversea chinese banking corp ltd raise position in ecolab inc nyseecl westchester capital management inc ha million stake in becton dickinson and company nysebdx ecolab inc nyseecl share purchased by mirabaud cie sa limited option for democrat to retaliate if texas republican redraw congressional map trump hit brazilian product with tariff over bolsonaro sign that someone you know is secretly broke do they apply to the people around you party poopers less than of american friendship cross political line texas house committee advance redrawn congressional map congressional redistricting map pass house committee pushing dems towards quorum break how trump became the new master of the senate luka doncic lakers agree to year million extension banco santander sa grows stock position in the walt disney company nysedis union pacific corporation nyseunp share sold by mitsubishi ufj trust banking corp moor cabot inc ha million stock holding in zoetis inc nysezts black swimmer teach others to help end history of aquatic segregation laduke let the sun shine so say former intelligence agent rolling stone claim the kremlin made up russiagate evidence can donald trump f diplomacy keep everyone happy laduke let the sun shine laduke let the sun shine in japan launched a second strike on pearl harborsort of laduke let the sun shine your view people have the right to sign recall petition laduke let the sun shine letter to the editor the generosity of kentuckian laduke let the sun shine laduke let the sun shine it super disrespectful wnba player exasperated after sex toy thrown onto court for second time in a week it trump economy now the latest financial number offer some warning sign aviso wealth management sell share of eli lilly and company nyselly j safra sarasin holding ag buy share of procter gamble company the nysepg steve witkoff meet with hostage family in tel aviv u court upholds order blocking indiscriminate targeting by immigration patrol u court upholds order blocking indiscriminate targeting by immigration patrol oversea chinese banking corp ltd reduces stock position in welltower inc nysewell here what we know about the anaconda shooter the victim and ongoing search effort banque cantonale vaudoise purchase new position in entergy corporation nyseetr the sherwinwilliams company nyseshw position lifted by banco santander sa lsv asset management sell share of phillips nysepsx profit drop at warren buffetts berkshire hathaway a it writes down it kraft heinz investment moor cabot inc sell share of thermo fisher scientific inc nysetmo trump tariff are making money that may make them hard to quit the new york time revolution camp unhinged white liberal spend weekend smashing junk in antitrump fury sen josh hawleys longshot stock trading ban bill advance rankles trump larimer county restaurant inspection made july st louis boeing worker again reject contract strike to begin at midnight former associate executive director file lawsuit against ahsaa boeings st louis union worker will strike on monday boeing defense union poised for it first strike since bloombergcom mahjong the latest networking opp for woman entrepreneur why analyst say this audited ai token is the next big thing coinmarketcap listing spark daily sale in the million luxury label lose their luster amid changing vibe gucci sale plunge and lvmh disappoints best crypto to buy now and secure early gain north texas father of win million at choctaw casino this cookware maker is bracing for steel tariff behind a wall of pot and pan nbc news opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day trump live update white house defends firing bls commissioner over weak job number the new york time opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day white house official defend firing of labor official a critic warn of trust erosion clickbait and clean water youtubers are the new philanthropist palantir advanced micro walt disney pfizer zoominfo and more stock to watch this week barrons sen blackburn to introduce bill to root out embedded foreign interest free iq test online with instant result reliable accurate free realiqonline roll out new iq testing experience himerex pro this himerex pro uk app set new standard in aidriven trading with unmatched security and user approval ev lead vehicle category in california berkshire earnings key takeaway strong profit no buyback high cash level barrons steam and itchio pull hundred of adult game over alleged visa mastercard pressure pauline ferrandprevot win tour de france femmes becoming first french rider to win event since iam district member in st louis reject latest boeing offer strike for fair contract iam union bmw sue teen entrepreneur twice over alibaba trademark infringement top trader see ripple xrp hitting binance coin bnb and little pepe lilpepe when to buy for maximum roi starbucks to close pickup store in back to starbucks pivot citigroup launch stratum elite card to rival chase sapphire reserve disney world ferry crash into dock in stormy weather injuring rider u manufacturing resurgence foreign firm relocate amid tariff the fantastic four first step hold it lead atop the box office the fantastic four first step hold it lead atop the box office the fantastic four first step hold it lead atop the box office la vega unemployment at amid tourism drop and visitor de'
$ws.Cells.Item(29, 10).Value = 'nintendo break sale record delayed gta release save switch fedex corporation fdxs ceo hitter jim cramer inside gamesxs ai revolution india become next global gamechanger intel corporation intc taught trust ai chatbots jim cramer ipo stock china stock atour test buy amid market fall procter gamble company pg great jim cramer espn obtains redzone medium asset nfl blockbuster deal per report ai agentic ai adoption plus ai action plan china mcdonalds plan major ai expansion streamline customer order executive paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame oil major beat billion fine kazakhstan paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame outdated mindset business paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame worked building people killed building secure paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame columbia selling stylish quilted jacket reviewer wellmade paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame action comedy anime series paula deen abruptly closed savannah restaurant launched network fame paula deen abruptly closed savannah restaurant launched network fame hit k limit next cadence design system cdnss nvidia great shape china jim cramer alnylams amvuttra powerful entrance eye firstline dominance president trump seek firing government overseeing job data weak employment report president trump seek firing government overseeing job data weak employment report nvidia corporation nvda great shape china jim cramer boeing company ba horse ceo great job jim cramer doj identifies tariff fraud enforcement priority waste wms ceo remarkable job jim cramer merck co mrk move talking animal division jim cramer starbucks corporation sbux dont worry luckin coffee jim cramer sarah jessica parker announces emotional poetic end carrier global corporation carr delivered good quarter jim cramer iphone pro launch september colour pro chip expected india price otis worldwide corporation otis tough quarter china horrible jim cramer apple watch series click away stacksocial celestica cl jim cramer ai evolution prompt sense urgency among cfo oracle jury tesla partly blame fatal crash trump issue nuclear threat russiaand wall street fear gauge vertical crypto founder take note sec roundtable decide rule govern sugar price pressured brazil ramp sugar production stablecoins building nextgen rail enterprise finance tesla must pay million damage fatal autopilot crash jury cnbc miami jury find tesla partly liable autopilot crash award million punitive damage nbc bond catch tailwind job data tesla tsla found liable fatal autopilot crash pay million latest stock market global trade partner react trump tariff trump fire labor statistic bos weak revised job data green revolution anand chandra cofounder executive director aryaag fixing india agri gap tech cocoa price sharply lower outlook adequate supply ghislaine maxwell move lowersecurity federal prison camp signal reddit replace google internet search engine investor ct bond commission approves lower electric bill republican enough coffee price plunge supply concern ease latest business trading partner worldwide respond trump tariff latest business trading partner worldwide respond trump tariff nextnav oscilloquartz partner advance g timing solution critical infrastructure resmed analyst increase forecast strong q earnings hedge fund insider trading william heard warren buffett ray dalio bill ackman rokos alden global channel outdoor cco buy stock friday nosediving market coinbase stock coin cratering tipranks phoenix az lemonade stand squeeze success expands western million sale cryptos join deliver bowling center west valley bowl face endofmonth deadline corporation public broadcasting shutter trump funding cut best trader congress track politician stock trade dallas bank confident loss pay best gold stock watch midtiers junior beat major online creator led mrbeast mark rober raise million clean water access online creator led mrbeast mark rober raise million clean water access online creator led mrbeast mark rober raise million clean water access chairman savannah economic development dy jury order tesla pay autopilot crash opening costly lawsuit jury order tesla pay autopilot crash opening costly lawsuit online creator led mrbeast mark rober raise million clean water access tesla found partially liable fatal crash involving driverassistance technology jury order million damage washington post india secret weapon defence sovereignty ai zerochina tech trump injects dose uncertainty tariff push start date aug trump injects dose uncertainty tariff push start date aug analyst boost forecast microstrategy following q trump injects dose uncertainty tariff push start date aug unit moving portable storage primary sponsor timmy hill watkins glen nascar truck series race tech stock lead selloff soft labor data fuel fed cut bet eye hurricane job market soured economist expensive trump massive tariff hike expensive trump massive tariff hike trump seek goldplated stain white house beatrice police department hosting national night nicole leone executive director erie humane society september trial date president trump seek firing government overseeing job data weak employment report president trump seek firing government overseeing job data weak employment report president trump seek firing government overseeing job data weak employment report president trump seek firing government overseeing job data weak employment report president trump seek firing government overseeing job data weak employment report president trump seek firing government overseeing job data weak employment report north carolina senate race turn battle middle class sarah jessica parker announces emotional poetic end ghislaine maxwell join famous prisoner move texas prison camp wave fake credential spark political fallout spain gop congressman suggests hed lose money broke adviser stop stock trade gop congressman suggests hed lose money broke adviser stop stock trade gop congressman suggests hed lose money broke adviser stop stock trade jeffrey epstein girlfriend ghislaine maxwell transferred prison camp texas jeffrey epstein girlfriend ghislaine maxwell transferred prison camp texas jeffrey epstein girlfriend ghislaine maxwell transferred prison camp texas jeffrey epstein girlfriend ghislaine maxwell transferred prison camp texas gop congressman suggests hed lose money broke adviser stop stock trade gop congressman suggests hed lose money broke adviser stop stock trade gop congressman suggests hed lose money broke adviser stop stock trade corporation public broadcasting fund pb npr close federal aid cut nyc subway nightmare flood trap commuter filth failure latest stock market global trade partner react trump tariff ghislaine maxwell move lowersecurity federal prison camp signal public rule speaking local meeting bismarckmandan area corporation public broadcasting shutter trump funding cut congressman view truth impact bill st louis county trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev college rally planned another local health system cut youth trans care american couple turning italy dolce vita quest memorable wedding thursday odor came construction industrial water plant trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia cult cruelty ghastly photo video emaciated israeli hostage released gaza terrorist expensive trump massive tariff hike red move slash medicaid funding planned parenthood trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev trump order nuclear sub repositioned statement exrussian medvedev guest columnist must empower support girl boy trump order nuclear submarine repositioned provocative statement russia editorial maxwell giving life conspiracy shribman hero age trump brown fulfilling promise ada trophy hunter killed lion zimbabwe part research project sparking anger pensacola man traveled orlando stealing smokey bear sign along kamala harris insufferably cackle government spying american wireless earpods western country speak future palestinian nightmare unfolding gaza worsens western country speak future palestinian nightmare unfolding gaza worsens western country speak future palestinian nightmare unfolding gaza worsens western country speak future palestinian nightmare unfolding gaza worsens western country speak future palestinian nightmare unfolding gaza worsens trump deploying nuclear sub provocative exrussian democrat struggling block trump texas power grab healing power sunlight office cleaner came facetoface manhattan gunman detail deadly rampage corporation public broadcasting funder npr pb end operation month federal budget cut dozen protester swarm chuck schumer kirsten gillibrands nyc office gaza fourth circuit rule unanimously favor employee involving employer allegation breach restrictive covenant county animal control texas researcher airport face deportation texas researcher airport face deportation texas researcher airport face deportation texas researcher airport face deportation texas researcher airport face deportation trump plan build million white house party room america rising part federal judge trump administration lying supreme k nieth retire fbi reportedly redacted reference trump epstein file trump positioning nuclear submarine near russia following russian threat war immigration legislation dangerous immigrant judge pause removal k migrant honduras nicaragua nepal overthetop ruling democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan democrat launching summer blitz press republican trump spending plan trump weighs american eagle sydney sweeney ad shiba inu cardano set breakout token set sight gain cryptos buy trimiq weight loss launched uk ireland green tea probiotic complex free spin deposit bonus casino online real money casino deposit wild casino shock job report set recession alert fresh clue ai boosting unemployment jpmorgan fortune crypto presale millionaire dogecoin ada xrp cant match ozak ai potential arialief neuropathy supplement launch sciatic nerve discomfort relief ig choice rail reader commentary youre risk wage garnished student loan debt raised climbing ruvi ai ruvi ripple xrp coinmarketcap listing token highest demand parent suing snapchat fentanyl death tea encouraged user spill apps data got leaked npr federal reserve governor resigns creating vacancy trump federal reserve governor resigns creating vacancy trump upicom bragar eagel squire pc reminds investor class action neogen alto bragar eagel squire pc reminds investor class action neogen alto rxsight replimune encourages investor contact firm audizen drop officially launched natural tinnitus support inspired audifort europe ev sale report tesla plus ev market matsato knife launch australia usa redefining precision power kitchen warren buffett dividend stock buy forever crushed q earnings luka doncic agrees threeyear million contract extension la lakers financial swamp banking regulator decide american bank account shut sam altman compare openai manhattan projectand joking risk democrat mamdani problem editorial toncoin ondo cross cold wallet entry crypto reward outclass weekend gaming best ipad game try little pepe lilpepe price prediction strong potential beat solana sol cardano ada sluggish summer sin la vega drop tourist business pm business pm lao brazil trump tariff lot loser winner pay price lao brazil trump tariff lot loser winner pay price lao brazil trump tariff lot loser winner pay price lao brazil trump tariff lot loser winner pay price texas panel advance redrawn congressional map take democratic seat money exit business equity without regret pepes chart turn positive litecoin gather whale support blockdag expands web reach devs workweek isnt good team moral good bottom line teleflex raise outlook following strong q biotronik acquisition tucson trolley tour operator seek expand reach trump prosecutor jack smith probed alleged election interference hiring slowed sharply summer wall street journal mstr michael saylor explains strategy strc preferred stock firm iphone moment coindesk harvard comply dhs request employment form missed binance coin bnb expert ruvi ai ruvi audited token next millionaire maker coinmarketcap listing spark high demand youre risk wage garnished student loan debt safely download install game laptop shib find support ondo cold wallet price gas rebate crypto watch little pepe lilpepe gather momentum gain ahead cardano ada nears run trump administration take orwellian turn avalanche washington wall street kept investor guessing westchester trim stake valero energy corporation nysevlo united super pty ltd capacity trustee construction building union superannuation fund investment valero energy corporation nysevlo decadeold bitcoin wallet emptied million btc moved amid market correction florida quietly surpassed california solar growth ukrainian longrange drone hit multiple russian refinery factory overnight journointerior designer jim acosta guarantee trump wh ballroom look st americu appoints community relation manager mitsubishi ufj trust banking corp purchase ge aerospace nysege j safra sarasin ag purchase thermo fisher scientific nysetmo mitsubishi ufj trust banking corp sell thermo fisher scientific nysetmo banco santander sa million stake thermo fisher scientific nysetmo walt disney company nysedis bought bank nova scotia trust co walt disney company nysedis stake lifted r squared ltd walt disney company nysedis acquired j safra sarasin ag walt disney company nysedis bought j safra sarasin ag banco santander sa increase stock walt disney company nysedis ai set flip major altcoinare figma ipo surprise winner charity million sharesand famous backstory sparked bitter feud oil fortune decade ago fortune hidden signal oil market moor cabot million stock goldman sachs nysegs moor cabot buy goldman sachs nysegs hartford investment co cut stock zoetis nysezts florida quietly surpassed california solar growth cnbc lesson decade empowering woman entrepreneur lead boldly best trucking bookkeeping service stock consider blockdags vesting pa roi end soon eth eye jup drop smallbusiness owner child care crisis critical problem pepeto price forecast presale meme coin surpass dogecoin shiba inu ethereum editorial beyond study subminimum must consider ramification palantir pltr report q earnings best trucking bookkeeping service super disrespectful wnba player exasperated sex toy thrown onto claude cutoff anthropic revoke openais api access ahead gpt launch malaga cove sell abbott laboratory nyseabt bayforest ltd stake fidelity national financial nysefnf vestcor buy fidelity national financial nysefnf modera wealth sell arista network nyseanet bank nova scotia sell procter gamble company nysepg bank nova scotia increase stock arista network nyseanet aviso wealth acquires procter gamble company nysepg aviso wealth cut boeing company nyseba j safra sarasin ag purchase abbott laboratory nyseabt pathway financial adviser purchase eli lilly company nyselly aviso wealth lower stock boeing company nyseba procter gamble company nysepg bought j safra sarasin ag best dividend growth stock invest hims alert bragar eagel squire pc announces class action lawsuit steve witkoff meet hostage family tel aviv susan estrich harvard fight money v principle jet holiday ad audio used deportation video shes cringe christine pelosi tease shes running youd video dream dr king nightmare donald trump opinion dream dr king nightmare donald trump opinion dream dr king nightmare donald trump opinion guam printing facility enhance navy maintenance capability ig choice rail reader commentary trump motif opaque feminist economist bare india stand firm russian oil import trump sanction threat pedaling dei cycling weekly demand bike community include fat black woman theyre racist military launched largest air exercise pennsylvania lawmaker child care broken dont agree fix yet twice tale council attends team building workshop trump save maxwell destroy everything promised pope thrill hundred thousand young catholic holy youth festival pope thrill hundred thousand young catholic holy youth festival republican enough confirm trump nominee conservative critic texas panel advance redrawn congressional map take democratic seat trump press secretary performance answer took weird turn corvallis nothing help gazans democrat gearing online texas moving forward radical redistricting plan america trucker place park george mason university president job raise trump admins probe bls biased donna brazile nprpbsdefunding maga compete idea cue hard trip hillary degree graduate face discouraging odds landing job trump administration move ban abortion care veteran unspeakably cruel manhunt intensifies suspect montana bar mass shooting texas panel advance redrawn congressional map take democratic seat bill maher trump losing credibility supporter smell bullst kevin oleary slam trump whacking bls chief disappointing job report dont shoot messenger trump administration halted civil lawsuit targeting abuse prisoner mentally ill never safety whats blocked uk orwellian online bill country dominate fighter jet export marketfor tiktok shape america next generation beijing trump prosecutor jack smith probed alleged election interference virginia end fiscal million surplus next budget bonus editorial cartoon aug voice conserving enhancing agricultural resource solar energy trump alum chilling google message secondterm lawfare finest texas house panel advance redrawn congressional map add gop seat nd representative helping put nation closer fiscal crisis tracey weiss ocean backyard coming soon opening gate cotonicoast dairy senate progress averting government shutdown usual piecemeal deal witkoff hostage family trump full gaza agreement wnba approve record deal connecticut sun boston proud democrat call chris cuomos admits trump shes texas house committee advance gopfriendly map texas house committee advance gopfriendly map texas house committee advance gopfriendly map release epstein file end coverup texas house committee advance gopfriendly map texas house committee advance gopfriendly map texas house committee advance gopfriendly map trump administration take orwellian turn photo reveal iconic white house rose garden paved trump makeover nfls kickoff rule stay slight tweak else changing oversea chinese banking corp ltd raise ecolab nyseecl westchester million stake becton dickinson company nysebdx ecolab nyseecl purchased mirabaud cie sa limited option democrat retaliate texas republican redraw congressional map trump hit brazilian product tariff bolsonaro sign someone secretly broke apply people around party poopers less american friendship cross political line texas house committee advance redrawn congressional map congressional redistricting map pas house committee pushing dems towards quorum break trump became master senate luka doncic lakers agree million extension banco santander sa grows stock walt disney company nysedis union pacific corporation nyseunp mitsubishi ufj trust banking corp moor cabot million stock zoetis nysezts black swimmer teach others help end history aquatic segregation laduke sun shine intelligence agent rolling stone claim kremlin russiagate evidence donald trump f diplomacy everyone happy laduke sun shine laduke sun shine japan launched strike pearl harborsort laduke sun shine view people sign recall petition laduke sun shine editor generosity kentuckian laduke sun shine laduke sun shine super disrespectful wnba player exasperated sex toy thrown onto trump economy latest financial offer warning sign aviso wealth sell eli lilly company nyselly j safra sarasin ag buy procter gamble company nysepg steve witkoff meet hostage family tel aviv upholds order blocking indiscriminate targeting immigration patrol upholds order blocking indiscriminate targeting immigration patrol oversea chinese banking corp ltd reduces stock welltower nysewell anaconda shooter victim ongoing search effort banque cantonale vaudoise purchase entergy corporation nyseetr sherwinwilliams company nyseshw lifted banco santander sa lsv asset sell phillips nysepsx profit drop warren buffetts berkshire hathaway writes kraft heinz investment moor cabot sell thermo fisher scientific nysetmo trump tariff making money hard quit york revolution camp unhinged white liberal spend weekend smashing junk antitrump fury sen josh hawleys longshot stock trading ban bill advance rankles trump larimer county restaurant inspection st louis boeing worker reject contract strike begin midnight associate executive director file lawsuit ahsaa boeings st louis union worker strike monday boeing defense union poised strike since bloombergcom mahjong latest networking opp woman entrepreneur analyst audited ai token next coinmarketcap listing spark sale million luxury label lose luster amid changing vibe gucci sale plunge lvmh disappoints best crypto buy secure early gain north texas father million choctaw casino cookware maker bracing steel tariff behind wall pot pan nbc opec country boost oil production barrel per opec country boost oil production barrel per opec country boost oil production barrel per opec country boost oil production barrel per opec country boost oil production barrel per trump update white house defends firing bls commissioner weak job york opec country boost oil production barrel per opec country boost oil production barrel per white house defend firing labor critic warn trust erosion clickbait clean water youtubers philanthropist palantir advanced micro walt disney pfizer zoominfo stock watch barrons sen blackburn introduce bill root embedded foreign interest free iq test online instant reliable accurate free realiqonline roll iq testing experience himerex pro himerex pro uk app set aidriven trading unmatched security user approval ev lead vehicle category california berkshire earnings takeaway strong profit buyback high cash level barrons steam itchio pull hundred adult game alleged visa mastercard pressure pauline ferrandprevot tour de france femmes becoming french rider event since iam district st louis reject latest boeing offer strike fair contract iam union bmw sue teen entrepreneur twice alibaba trademark infringement trader ripple xrp hitting binance coin bnb little pepe lilpepe buy maximum roi starbucks close pickup store starbucks pivot citigroup launch stratum elite card rival chase sapphire reserve disney world ferry crash dock stormy weather injuring rider manufacturing resurgence foreign firm relocate amid tariff fantastic step lead atop box office fantastic step lead atop box office fantastic step lead atop box office la vega unemployment amid tourism drop visitor de'
